# Apply weekly row-data permutation for columns D, I, J, K, L, M, N, O, P
# across data rows 2-46 (header is row 1). Columns A,B,C,E,F,G,H,Q,R are
# identical for every row and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","I","J","K","L","M","N","O","P")

# Mapping of destination row -> source row (values that were at $Src in the
# original sheet must end up at $Dest).
$mapping = @(
    @{ Dest = 2; Src = 37 },
    @{ Dest = 3; Src = 28 },
    @{ Dest = 4; Src = 39 },
    @{ Dest = 5; Src = 20 },
    @{ Dest = 6; Src = 21 },
    @{ Dest = 7; Src = 10 },
    @{ Dest = 8; Src = 2 },
    @{ Dest = 9; Src = 4 },
    @{ Dest = 10; Src = 23 },
    @{ Dest = 11; Src = 12 },
    @{ Dest = 12; Src = 15 },
    @{ Dest = 13; Src = 38 },
    @{ Dest = 14; Src = 43 },
    @{ Dest = 15; Src = 16 },
    @{ Dest = 16; Src = 17 },
    @{ Dest = 17; Src = 30 },
    @{ Dest = 18; Src = 26 },
    @{ Dest = 19; Src = 13 },
    @{ Dest = 20; Src = 14 },
    @{ Dest = 21; Src = 8 },
    @{ Dest = 22; Src = 27 },
    @{ Dest = 23; Src = 29 },
    @{ Dest = 24; Src = 33 },
    @{ Dest = 25; Src = 44 },
    @{ Dest = 26; Src = 9 },
    @{ Dest = 27; Src = 32 },
    @{ Dest = 28; Src = 22 },
    @{ Dest = 29; Src = 31 },
    @{ Dest = 30; Src = 11 },
    @{ Dest = 31; Src = 45 },
    @{ Dest = 32; Src = 5 },
    @{ Dest = 33; Src = 6 },
    @{ Dest = 34; Src = 7 },
    @{ Dest = 35; Src = 18 },
    @{ Dest = 36; Src = 19 },
    @{ Dest = 37; Src = 25 },
    @{ Dest = 38; Src = 34 },
    @{ Dest = 39; Src = 35 },
    @{ Dest = 40; Src = 36 },
    @{ Dest = 41; Src = 3 },
    @{ Dest = 42; Src = 46 },
    @{ Dest = 43; Src = 40 },
    @{ Dest = 44; Src = 41 },
    @{ Dest = 45; Src = 42 },
    @{ Dest = 46; Src = 24 }
)

# 1) Snapshot the original values of every affected cell before writing
#    anything (since the mapping is a permutation, a row used as a source
#    for one destination may itself be overwritten later).
$snapshot = @{}
foreach ($r in 2..46) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# 2) Write the permuted values using the snapshot so every destination gets
#    the original source row's data regardless of write order.
foreach ($entry in $mapping) {
    $dest = $entry.Dest
    $src = $entry.Src
    foreach ($c in $cols) {
        $srcAddr = "$c$src"
        $destAddr = "$c$dest"
        $ws.Range($destAddr).Value = $snapshot[$srcAddr]
    }
}
